$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 23:22"

# Update Cataluña row (row 5) statistics
$ws.Range("B5").Value = 28323
$ws.Range("C5").Value = 10378
$ws.Range("D5").Value = 15037
$ws.Range("E5").Value = 2908
